$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# " Monday, 20 September, 2021" -> " " + "Monday, 20 September, 2021"
# (new run gets identical explicit formatting: JetBrains Mono NL, bold,
# green 00BF00, sz 36 -- same as the run it was split from)
$find1 = $d.Content.Find
[void]$find1.Execute("Monday, 20 September, 2021", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$dateRng = $find1.Parent

# Toggle the color so the interop layer actually materializes a run
# split (setting it straight to the already-applied value is a no-op).
$dateRng.Font.Color = 1
$dateRng.Font.Color = 48896

# --- Change 2 -----------------------------------------------------------
# " is located." -> " is located" + "."
# (new run gets identical explicit formatting: JetBrains Mono NL,
# italic, red FF0000, sz 18 -- same as the run it was split from)
$find2 = $d.Content.Find
[void]$find2.Execute(" is located.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$locatedRng = $find2.Parent
$dotRng = $d.Range($locatedRng.End - 1, $locatedRng.End)

$dotRng.Font.Color = 1
$dotRng.Font.Color = 255
